$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (shared string) into a cell that is
# subsequently styled with the numeric "0.00000" format, matching how
# the source workbook stores these benchmark numbers (text-typed, but
# carrying a numeric display format).
function Set-NumFmtTextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value2 = $val
    $r.NumberFormat = "0.00000"
}

# Plain text cell with no special number format (used for the brand new row 14).
function Set-PlainTextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value2 = $val
    $r.ClearFormats()
}

# ---- Row 11 (threads = 100) ----
Set-NumFmtTextCell "C11" "0.000192"
Set-NumFmtTextCell "D11" "0.000426"
Set-NumFmtTextCell "E11" "0.000763"
Set-NumFmtTextCell "F11" "0.00108"

# ---- Row 12 (threads = 1000) ----
Set-NumFmtTextCell "C12" "0.00061"
Set-NumFmtTextCell "D12" "0.000607"
Set-NumFmtTextCell "E12" "0.000800"
Set-NumFmtTextCell "F12" "0.00133"

# ---- Row 13 (threads = 10000) ----
Set-NumFmtTextCell "C13" "0.00349"
Set-NumFmtTextCell "D13" "0.00289"
Set-NumFmtTextCell "E13" "0.00243"
Set-NumFmtTextCell "F13" "0.00196"

# ---- Row 14 (threads = 100000) - brand new row, default style ----
$ws.Range("B14").Value2 = 100000
Set-PlainTextCell "C14" "0.02699"
Set-PlainTextCell "D14" "0.0171"
Set-PlainTextCell "E14" "0.01313"
Set-PlainTextCell "F14" "0.01166"

# ---- H8: "runs: 1000" -> "runs: 100" ----
# (set last so the new shared string lands at the end of sharedStrings.xml)
$ws.Range("H8").Value2 = "runs: 100"

# ---- Column widths for the new data columns ----
$ws.Range("C1:F1").EntireColumn.AutoFit()

# ---- View: active cell / selection moves to H11 ----
$ws.Range("H11").Select()

# ---- Zoom the sheet to 80% ----
$wb.Windows.Item(1).Zoom = 80

# ---- Page setup: paper size / orientation ----
$ws.PageSetup.PaperSize = 9      # xlPaperA4
$ws.PageSetup.Orientation = 1    # xlPortrait

Write-Output "done"
